$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel scenario names in column A (rows 2-7): shift bN -> b(N-1)
$ws.Range("A2").Value = "b0"
$ws.Range("A3").Value = "b1"
$ws.Range("A4").Value = "b2"
$ws.Range("A5").Value = "b3"
$ws.Range("A6").Value = "b4"
$ws.Range("A7").Value = "b5"
